$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.729.74"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "'2.535.22"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'565.14"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").Value = "'150.22"
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.611"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").Value = "'2.533.90"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").Value = "'0.112"
$ws.Range("E10").Value = "  -5.75%  "
$ws.Range("D11").Value = "'5.52"
$ws.Range("E11").Value = "  -5.00%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "'0.369"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").Value = "'26.96"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "'2.991.22"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "'0.0000172"
$ws.Range("E16").Value = "  -5.00%  "
$ws.Range("D17").Value = "'62.534.83"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "'2.522.79"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").Value = "'11.50"
$ws.Range("E19").Value = "  -4.41%  "
$ws.Range("D20").Value = "'7.23"
$ws.Range("E20").Value = "  -5.09%  "
$ws.Range("E21").Value = "  -4.36%  "
$ws.Range("D22").Value = "'329.22"
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("D24").Value = "'65.50"
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("D25").Value = "'1.82"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("D26").Value = "'0.0000108"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "'1.60"
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("D28").Value = "'2.648.41"
$ws.Range("D29").Value = "'8.77"
$ws.Range("E29").Value = "  -4.34%  "
$ws.Range("D30").Value = "'558.80"
$ws.Range("E30").Value = "  -5.10%  "
$ws.Range("D31").Value = "'8.02"
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").Value = "'0.156"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").Value = "'1.96"
$ws.Range("E34").Value = "  -5.00%  "
$ws.Range("E35").Value = "  -5.56%  "
$ws.Range("D36").Value = "'6.16"
$ws.Range("E36").Value = "  -5.49%  "
$ws.Range("D37").Value = "'5.05"
$ws.Range("E37").Value = "  -5.48%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "'0.390"
$ws.Range("E39").Value = "  -3.92%  "
$ws.Range("D40").Value = "'18.97"
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("D41").Value = "'151.77"
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("D42").Value = "'1.78"
$ws.Range("E42").Value = "  -4.58%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'41.05"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").Value = "'2.40"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").Value = "'152.95"
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").Value = "'3.75"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("D48").Value = "'22.51"
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("D49").Value = "'0.0559"
$ws.Range("E49").Value = "  -5.05%  "
$ws.Range("D50").Value = "'0.607"
$ws.Range("E50").Value = "  -3.14%  "
$ws.Range("D51").Value = "'0.0966"
$ws.Range("E51").Value = "  -4.70%  "
